$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns A-D, rows 2-11
$values = @(
    @(2, 1, 5, 5),
    @(4, 1, 10, 10),
    @(5, 1, 15, 15),
    @(10, 3, 5, 5),
    @(6, 5, 5, 5),
    @(8, 5, 10, 10),
    @(9, 5, 15, 16),
    @(3, 6, 5, 5),
    @(1, 7, 5, 5),
    @(7, 7, 10, 10)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowValues = $values[$i]
    $ws.Cells.Item($row, 1).Value = $rowValues[0]
    $ws.Cells.Item($row, 2).Value = $rowValues[1]
    $ws.Cells.Item($row, 3).Value = $rowValues[2]
    $ws.Cells.Item($row, 4).Value = $rowValues[3]
}
